# Updated cryptos list with GitHub Actions
# Refresh the Price (D) / Volume(1h) (E) columns with newly scraped values,
# and swap the PEPE / InternetComputer(DFINITY) rows (30-31) which changed
# rank order. Price-looking numeric strings are written with a leading
# apostrophe so Excel keeps them as text (preserving trailing zeros / the
# "thousands-dot" formatting) instead of coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.475.87"
$ws.Range("E2").Value = "  -2.07%  "
$ws.Range("D3").Value = "2.480.89"
$ws.Range("E3").Value = "  -3.22%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'566.03"
$ws.Range("E5").Value = "  -3.12%  "
$ws.Range("D6").Value = "'164.40"
$ws.Range("E6").Value = "  -4.60%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.510"
$ws.Range("E8").Value = "  -2.00%  "
$ws.Range("D9").Value = "2.482.04"
$ws.Range("E9").Value = "  -3.03%  "
$ws.Range("D10").Value = "'0.156"
$ws.Range("E10").Value = "  -6.04%  "
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Value = "'0.352"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("D13").Value = "'4.89"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "2.939.41"
$ws.Range("E14").Value = "  -2.94%  "
$ws.Range("D15").Value = "69.452.61"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("E16").Value = "  -3.14%  "
$ws.Range("D17").Value = "'24.14"
$ws.Range("E17").Value = "  -5.55%  "
$ws.Range("D18").Value = "2.487.65"
$ws.Range("E18").Value = "  -2.96%  "
$ws.Range("D19").Value = "'11.10"
$ws.Range("E19").Value = "  -3.38%  "
$ws.Range("D20").Value = "'7.33"
$ws.Range("E20").Value = "  -7.92%  "
$ws.Range("D21").Value = "'344.29"
$ws.Range("E21").Value = "  -3.64%  "
$ws.Range("D22").Value = "'3.83"
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("D23").Value = "'1.91"
$ws.Range("E23").Value = "  -7.87%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'70.09"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").Value = "'3.85"
$ws.Range("E26").Value = "  -5.60%  "
$ws.Range("D28").Value = "'8.60"
$ws.Range("E28").Value = "  -6.55%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'7.70"
$ws.Range("E30").Value = "  -3.54%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0863"
$ws.Range("E31").Value = "  -7.07%  "
$ws.Range("D32").Value = "'439.89"
$ws.Range("E32").Value = "  -7.60%  "
$ws.Range("D33").Value = "'1.18"
$ws.Range("E33").Value = "  -8.56%  "
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("E35").Value = "  -4.64%  "
$ws.Range("D36").Value = "'155.78"
$ws.Range("E36").Value = "  -2.41%  "
$ws.Range("E37").Value = "  -5.99%  "
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").Value = "'18.09"
$ws.Range("E39").Value = "  -3.56%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "'0.312"
$ws.Range("E41").Value = "  -3.56%  "
$ws.Range("D42").Value = "'1.57"
$ws.Range("E42").Value = "  -3.41%  "
$ws.Range("D43").Value = "'4.54"
$ws.Range("E43").Value = "  -7.29%  "
$ws.Range("D44").Value = "'37.92"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("D45").Value = "'2.14"
$ws.Range("E45").Value = "  -10.31%  "
$ws.Range("D46").Value = "'1.06"
$ws.Range("E46").Value = "  -9.44%  "
$ws.Range("D47").Value = "'139.03"
$ws.Range("E47").Value = "  -4.70%  "
$ws.Range("D48").Value = "'3.42"
$ws.Range("E48").Value = "  -3.92%  "
$ws.Range("D49").Value = "'0.510"
$ws.Range("E49").Value = "  -5.89%  "
$ws.Range("D50").Value = "'0.0727"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("D51").Value = "'0.571"
$ws.Range("E51").Value = "  -2.97%  "
